$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-12-23 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-24 Sunday", 2) | Out-Null
$d.Content.Find.Execute("7+91=98", $true, $false, $false, $false, $false, $true, 1, $false, "9-5=4", 2) | Out-Null
$d.Content.Find.Execute("15+55=70", $true, $false, $false, $false, $false, $true, 1, $false, "88-71=17", 2) | Out-Null
$d.Content.Find.Execute("86-21=65", $true, $false, $false, $false, $false, $true, 1, $false, "83+9=92", 2) | Out-Null
$d.Content.Find.Execute("45+0=45", $true, $false, $false, $false, $false, $true, 1, $false, "93-13=80", 2) | Out-Null
$d.Content.Find.Execute("15+68=83", $true, $false, $false, $false, $false, $true, 1, $false, "32-11=21", 2) | Out-Null
$d.Content.Find.Execute("85-64=21", $true, $false, $false, $false, $false, $true, 1, $false, "84-18=66", 2) | Out-Null
$d.Content.Find.Execute("57-18=39", $true, $false, $false, $false, $false, $true, 1, $false, "15-5=10", 2) | Out-Null
$d.Content.Find.Execute("57+3=60", $true, $false, $false, $false, $false, $true, 1, $false, "17+29=46", 2) | Out-Null
$d.Content.Find.Execute("8+31=39", $true, $false, $false, $false, $false, $true, 1, $false, "47-20=27", 2) | Out-Null
$d.Content.Find.Execute("60-29=31", $true, $false, $false, $false, $false, $true, 1, $false, "9+39=48", 2) | Out-Null
$d.Content.Find.Execute("49+45=94", $true, $false, $false, $false, $false, $true, 1, $false, "98-24=74", 2) | Out-Null
$d.Content.Find.Execute("58+12=70", $true, $false, $false, $false, $false, $true, 1, $false, "87-34=53", 2) | Out-Null
$d.Content.Find.Execute("31-9=22", $true, $false, $false, $false, $false, $true, 1, $false, "43+52=95", 2) | Out-Null
$d.Content.Find.Execute("61+10=71", $true, $false, $false, $false, $false, $true, 1, $false, "28+68=96", 2) | Out-Null
$d.Content.Find.Execute("22+10=32", $true, $false, $false, $false, $false, $true, 1, $false, "65+11=76", 2) | Out-Null
$d.Content.Find.Execute("22+22=44", $true, $false, $false, $false, $false, $true, 1, $false, "27-17=10", 2) | Out-Null
$d.Content.Find.Execute("2+18=20", $true, $false, $false, $false, $false, $true, 1, $false, "74+16=90", 2) | Out-Null
$d.Content.Find.Execute("57-50=7", $true, $false, $false, $false, $false, $true, 1, $false, "83+14=97", 2) | Out-Null
$d.Content.Find.Execute("23+2=25", $true, $false, $false, $false, $false, $true, 1, $false, "33+1=34", 2) | Out-Null
$d.Content.Find.Execute("34-14=20", $true, $false, $false, $false, $false, $true, 1, $false, "61+36=97", 2) | Out-Null
$d.Content.Find.Execute("61+30=91", $true, $false, $false, $false, $false, $true, 1, $false, "4+92=96", 2) | Out-Null
$d.Content.Find.Execute("71-30=41", $true, $false, $false, $false, $false, $true, 1, $false, "24-21=3", 2) | Out-Null
$d.Content.Find.Execute("53+5=58", $true, $false, $false, $false, $false, $true, 1, $false, "31+62=93", 2) | Out-Null
$d.Content.Find.Execute("97-56=41", $true, $false, $false, $false, $false, $true, 1, $false, "4+95=99", 2) | Out-Null
$d.Content.Find.Execute("9-0=9", $true, $false, $false, $false, $false, $true, 1, $false, "28+69=97", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $true, $false, $false, $false, $false, $true, 1, $false, "87-7=80", 2) | Out-Null
$d.Content.Find.Execute("36+63=99", $true, $false, $false, $false, $false, $true, 1, $false, "96-23=73", 2) | Out-Null
$d.Content.Find.Execute("52-42=10", $true, $false, $false, $false, $false, $true, 1, $false, "80-21=59", 2) | Out-Null
$d.Content.Find.Execute("79-27=52", $true, $false, $false, $false, $false, $true, 1, $false, "22+20=42", 2) | Out-Null
$d.Content.Find.Execute("17-7=10", $true, $false, $false, $false, $false, $true, 1, $false, "93-82=11", 2) | Out-Null
$d.Content.Find.Execute("86-81=5", $true, $false, $false, $false, $false, $true, 1, $false, "97-6=91", 2) | Out-Null
$d.Content.Find.Execute("54+24=78", $true, $false, $false, $false, $false, $true, 1, $false, "37+6=43", 2) | Out-Null
$d.Content.Find.Execute("91-42=49", $true, $false, $false, $false, $false, $true, 1, $false, "67-51=16", 2) | Out-Null
$d.Content.Find.Execute("53+15=68", $true, $false, $false, $false, $false, $true, 1, $false, "73+26=99", 2) | Out-Null
$d.Content.Find.Execute("24+74=98", $true, $false, $false, $false, $false, $true, 1, $false, "67-5=62", 2) | Out-Null
$d.Content.Find.Execute("68+17=85", $true, $false, $false, $false, $false, $true, 1, $false, "75-32=43", 2) | Out-Null
$d.Content.Find.Execute("99-61=38", $true, $false, $false, $false, $false, $true, 1, $false, "56-44=12", 2) | Out-Null
$d.Content.Find.Execute("3+8=11", $true, $false, $false, $false, $false, $true, 1, $false, "52+25=77", 2) | Out-Null
$d.Content.Find.Execute("39+56=95", $true, $false, $false, $false, $false, $true, 1, $false, "25+66=91", 2) | Out-Null
$d.Content.Find.Execute("34-9=25", $true, $false, $false, $false, $false, $true, 1, $false, "62-6=56", 2) | Out-Null
$d.Content.Find.Execute("11+52=63", $true, $false, $false, $false, $false, $true, 1, $false, "66-34=32", 2) | Out-Null
$d.Content.Find.Execute("2+51=53", $true, $false, $false, $false, $false, $true, 1, $false, "60-12=48", 2) | Out-Null
$d.Content.Find.Execute("60+5=65", $true, $false, $false, $false, $false, $true, 1, $false, "51-36=15", 2) | Out-Null
$d.Content.Find.Execute("58-2=56", $true, $false, $false, $false, $false, $true, 1, $false, "99-9=90", 2) | Out-Null
$d.Content.Find.Execute("17-15=2", $true, $false, $false, $false, $false, $true, 1, $false, "47+11=58", 2) | Out-Null
$d.Content.Find.Execute("13+14=27", $true, $false, $false, $false, $false, $true, 1, $false, "5+57=62", 2) | Out-Null
$d.Content.Find.Execute("51+6=57", $true, $false, $false, $false, $false, $true, 1, $false, "18+25=43", 2) | Out-Null
$d.Content.Find.Execute("97-27=70", $true, $false, $false, $false, $false, $true, 1, $false, "99-3=96", 2) | Out-Null
$d.Content.Find.Execute("71-48=23", $true, $false, $false, $false, $false, $true, 1, $false, "65-51=14", 2) | Out-Null
$d.Content.Find.Execute("79-36=43", $true, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2) | Out-Null
$d.Content.Find.Execute("5+29=34", $true, $false, $false, $false, $false, $true, 1, $false, "52+21=73", 2) | Out-Null
$d.Content.Find.Execute("64-15=49", $true, $false, $false, $false, $false, $true, 1, $false, "90-72=18", 2) | Out-Null
$d.Content.Find.Execute("26-20=6", $true, $false, $false, $false, $false, $true, 1, $false, "83-65=18", 2) | Out-Null
$d.Content.Find.Execute("89-87=2", $true, $false, $false, $false, $false, $true, 1, $false, "90-45=45", 2) | Out-Null
$d.Content.Find.Execute("34+15=49", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=7", 2) | Out-Null
$d.Content.Find.Execute("96-13=83", $true, $false, $false, $false, $false, $true, 1, $false, "42+12=54", 2) | Out-Null
$d.Content.Find.Execute("60-15=45", $true, $false, $false, $false, $false, $true, 1, $false, "30-2=28", 2) | Out-Null
$d.Content.Find.Execute("64-58=6", $true, $false, $false, $false, $false, $true, 1, $false, "35+11=46", 2) | Out-Null
$d.Content.Find.Execute("52+12=64", $true, $false, $false, $false, $false, $true, 1, $false, "64-13=51", 2) | Out-Null
$d.Content.Find.Execute("21+33=54", $true, $false, $false, $false, $false, $true, 1, $false, "72-16=56", 2) | Out-Null
$d.Content.Find.Execute("49+35=84", $true, $false, $false, $false, $false, $true, 1, $false, "47+48=95", 2) | Out-Null
$d.Content.Find.Execute("97-69=28", $true, $false, $false, $false, $false, $true, 1, $false, "48-17=31", 2) | Out-Null
$d.Content.Find.Execute("44-24=20", $true, $false, $false, $false, $false, $true, 1, $false, "68+21=89", 2) | Out-Null
$d.Content.Find.Execute("14+29=43", $true, $false, $false, $false, $false, $true, 1, $false, "92-83=9", 2) | Out-Null
$d.Content.Find.Execute("23+43=66", $true, $false, $false, $false, $false, $true, 1, $false, "88-79=9", 2) | Out-Null
$d.Content.Find.Execute("48-7=41", $true, $false, $false, $false, $false, $true, 1, $false, "61+15=76", 2) | Out-Null
$d.Content.Find.Execute("42-3=39", $true, $false, $false, $false, $false, $true, 1, $false, "64-16=48", 2) | Out-Null
$d.Content.Find.Execute("51-6=45", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=67", 2) | Out-Null
$d.Content.Find.Execute("38-23=15", $true, $false, $false, $false, $false, $true, 1, $false, "3+45=48", 2) | Out-Null
$d.Content.Find.Execute("4+1=5", $true, $false, $false, $false, $false, $true, 1, $false, "17-3=14", 2) | Out-Null
$d.Content.Find.Execute("11+58=69", $true, $false, $false, $false, $false, $true, 1, $false, "81-55=26", 2) | Out-Null
$d.Content.Find.Execute("49+48=97", $true, $false, $false, $false, $false, $true, 1, $false, "61+27=88", 2) | Out-Null
$d.Content.Find.Execute("14+26=40", $true, $false, $false, $false, $false, $true, 1, $false, "36-17=19", 2) | Out-Null
$d.Content.Find.Execute("45+6=51", $true, $false, $false, $false, $false, $true, 1, $false, "33+51=84", 2) | Out-Null
$d.Content.Find.Execute("57+24=81", $true, $false, $false, $false, $false, $true, 1, $false, "33-11=22", 2) | Out-Null
$d.Content.Find.Execute("62-46=16", $true, $false, $false, $false, $false, $true, 1, $false, "4+78=82", 2) | Out-Null
$d.Content.Find.Execute("82-56=26", $true, $false, $false, $false, $false, $true, 1, $false, "38+21=59", 2) | Out-Null
$d.Content.Find.Execute("86-40=46", $true, $false, $false, $false, $false, $true, 1, $false, "84+10=94", 2) | Out-Null
$d.Content.Find.Execute("48+51=99", $true, $false, $false, $false, $false, $true, 1, $false, "68-54=14", 2) | Out-Null
$d.Content.Find.Execute("39-1=38", $true, $false, $false, $false, $false, $true, 1, $false, "55-24=31", 2) | Out-Null
$d.Content.Find.Execute("2-2=0", $true, $false, $false, $false, $false, $true, 1, $false, "90-14=76", 2) | Out-Null
$d.Content.Find.Execute("7+37=44", $true, $false, $false, $false, $false, $true, 1, $false, "42+48=90", 2) | Out-Null
$d.Content.Find.Execute("39+4=43", $true, $false, $false, $false, $false, $true, 1, $false, "64+34=98", 2) | Out-Null
$d.Content.Find.Execute("0+62=62", $true, $false, $false, $false, $false, $true, 1, $false, "44-40=4", 2) | Out-Null
$d.Content.Find.Execute("45+1=46", $true, $false, $false, $false, $false, $true, 1, $false, "6+37=43", 2) | Out-Null
$d.Content.Find.Execute("62+30=92", $true, $false, $false, $false, $false, $true, 1, $false, "80-37=43", 2) | Out-Null
$d.Content.Find.Execute("25-4=21", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=41", 2) | Out-Null
$d.Content.Find.Execute("47+41=88", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=25", 2) | Out-Null
$d.Content.Find.Execute("57-46=11", $true, $false, $false, $false, $false, $true, 1, $false, "27-16=11", 2) | Out-Null
$d.Content.Find.Execute("30-8=22", $true, $false, $false, $false, $false, $true, 1, $false, "53+3=56", 2) | Out-Null
$d.Content.Find.Execute("36+38=74", $true, $false, $false, $false, $false, $true, 1, $false, "95-27=68", 2) | Out-Null
$d.Content.Find.Execute("93-86=7", $true, $false, $false, $false, $false, $true, 1, $false, "61+17=78", 2) | Out-Null
$d.Content.Find.Execute("69-46=23", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=92", 2) | Out-Null
$d.Content.Find.Execute("14+35=49", $true, $false, $false, $false, $false, $true, 1, $false, "26-10=16", 2) | Out-Null
$d.Content.Find.Execute("81-1=80", $true, $false, $false, $false, $false, $true, 1, $false, "7+92=99", 2) | Out-Null
$d.Content.Find.Execute("98-52=46", $true, $false, $false, $false, $false, $true, 1, $false, "31+51=82", 2) | Out-Null
$d.Content.Find.Execute("8+48=56", $true, $false, $false, $false, $false, $true, 1, $false, "73-20=53", 2) | Out-Null
$d.Content.Find.Execute("21+46=67", $true, $false, $false, $false, $false, $true, 1, $false, "5+5=10", 2) | Out-Null
$d.Content.Find.Execute("74+18=92", $true, $false, $false, $false, $false, $true, 1, $false, "17+47=64", 2) | Out-Null
$d.Content.Find.Execute("2+23=25", $true, $false, $false, $false, $false, $true, 1, $false, "42+35=77", 2) | Out-Null
